$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3253.25
$ws.Range("J88").Value = 3253.25
$ws.Range("L88").Value = 3253.25
$ws.Range("N88").Value = -4065.25
$ws.Range("H91").Value = 3253.25
$ws.Range("J91").Value = 3253.25
$ws.Range("L91").Value = 3253.25
$ws.Range("N91").Value = -6061.25
$ws.Range("H92").Value = 738.625
$ws.Range("I92").Value = 697.8276
$ws.Range("J92").Value = 1133
$ws.Range("K92").Value = 697.8276
$ws.Range("L92").Value = 1133
$ws.Range("M92").Value = 550.1724
$ws.Range("N92").Value = -3629
$ws.Range("H100").Value = 4443
$ws.Range("J100").Value = 5998
$ws.Range("L100").Value = 5998
$ws.Range("N100").Value = -7080
$ws.Range("H132").Value = 5482.0835
$ws.Range("I132").Value = 5546.5654
$ws.Range("K132").Value = 16639.6962
$ws.Range("M132").Value = -14109.6962
$ws.Range("H137").Value = 22600
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 22600
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 67800
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -72900

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2474.05
$ws.Range("I2").Value = 2041.8462
$ws.Range("J2").Value = 3276.7144
$ws.Range("K2").Value = 2041.8462
$ws.Range("L2").Value = 3276.7144
$ws.Range("M2").Value = -1928.8462
$ws.Range("N2").Value = -3502.7144
$ws.Range("H32").Value = 5648.5854
$ws.Range("I32").Value = 5160.7427
$ws.Range("K32").Value = 5160.7427
$ws.Range("M32").Value = -4873.7427
$ws.Range("H63").Value = 2118.8
$ws.Range("I63").Value = 1398.5
$ws.Range("K63").Value = 1398.5
$ws.Range("M63").Value = -712.5
$ws.Range("H66").Value = 2118.8
$ws.Range("I66").Value = 1398.5
$ws.Range("K66").Value = 6992.5
$ws.Range("M66").Value = -3560.5
$ws.Range("H116").Value = 2474.05
$ws.Range("I116").Value = 2041.8462
$ws.Range("J116").Value = 3276.7144
$ws.Range("K116").Value = 2041.8462
$ws.Range("L116").Value = 3276.7144
$ws.Range("M116").Value = 252.1538
$ws.Range("N116").Value = -7864.7144
$ws.Range("H132").Value = 29491928
$ws.Range("I132").Value = 13112.913
$ws.Range("J132").Value = 91129450
$ws.Range("K132").Value = 39338.739
$ws.Range("L132").Value = 273388350
$ws.Range("M132").Value = -36808.739
$ws.Range("N132").Value = -273393410

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2474.05
$ws.Range("I3").Value = 2041.8462
$ws.Range("J3").Value = 3276.7144
$ws.Range("K3").Value = 2041.8462
$ws.Range("L3").Value = 3276.7144
$ws.Range("M3").Value = -1927.8462
$ws.Range("N3").Value = -3504.7144
$ws.Range("H20").Value = 3142.25
$ws.Range("I20").Value = 1966
$ws.Range("J20").Value = 3534.3333
$ws.Range("K20").Value = 1966
$ws.Range("L20").Value = 3534.3333
$ws.Range("M20").Value = -1719
$ws.Range("N20").Value = -4028.3333
$ws.Range("H86").Value = 8382.32
$ws.Range("I86").Value = 11226
$ws.Range("K86").Value = 11226
$ws.Range("M86").Value = -10103
$ws.Range("H89").Value = 8382.32
$ws.Range("I89").Value = 11226
$ws.Range("K89").Value = 56130
$ws.Range("M89").Value = -50514
$ws.Range("H134").Value = 2595.2593
$ws.Range("I134").Value = 3115.9048
$ws.Range("K134").Value = 9347.714399999999
$ws.Range("M134").Value = -6812.714399999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 75000
$ws.Range("I26").Value = 75000
$ws.Range("K26").Value = 75000
$ws.Range("M26").Value = -74713
$ws.Range("H31").Value = 13892988
$ws.Range("I31").Value = 3051.4666
$ws.Range("J31").Value = 83342670
$ws.Range("K31").Value = 3051.4666
$ws.Range("L31").Value = 83342670
$ws.Range("M31").Value = -2756.4666
$ws.Range("N31").Value = -83343260
$ws.Range("H34").Value = 13892988
$ws.Range("I34").Value = 3051.4666
$ws.Range("J34").Value = 83342670
$ws.Range("K34").Value = 3051.4666
$ws.Range("L34").Value = 83342670
$ws.Range("M34").Value = -2849.4666
$ws.Range("N34").Value = -83343074
$ws.Range("H86").Value = 3997
$ws.Range("I86").Value = 3995
$ws.Range("J86").Value = 3998
$ws.Range("K86").Value = 3995
$ws.Range("L86").Value = 3998
$ws.Range("M86").Value = -2872
$ws.Range("N86").Value = -6244
$ws.Range("H89").Value = 3997
$ws.Range("I89").Value = 3995
$ws.Range("J89").Value = 3998
$ws.Range("K89").Value = 19975
$ws.Range("L89").Value = 19990
$ws.Range("M89").Value = -14359
$ws.Range("N89").Value = -31222
$ws.Range("H131").Value = 27238.8
$ws.Range("J131").Value = 27238.8
$ws.Range("L131").Value = 27238.8
$ws.Range("N131").Value = -37318.8
$ws.Range("H132").Value = 86291.71000000001
$ws.Range("J132").Value = 2292.4546
$ws.Range("L132").Value = 6877.3638
$ws.Range("N132").Value = -11937.3638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 21021.857
$ws.Range("I56").Value = 21021.857
$ws.Range("K56").Value = 21021.857
$ws.Range("M56").Value = -20491.857
$ws.Range("H86").Value = 779.7
$ws.Range("J86").Value = 1493.5
$ws.Range("L86").Value = 4480.5
$ws.Range("N86").Value = -6852.5
$ws.Range("H89").Value = 779.7
$ws.Range("J89").Value = 1493.5
$ws.Range("L89").Value = 13441.5
$ws.Range("N89").Value = -25297.5
$ws.Range("H134").Value = 4027.4285
$ws.Range("J134").Value = 11516.5
$ws.Range("L134").Value = 34549.5
$ws.Range("N134").Value = -44689.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 38463292
$ws.Range("I122").Value = 1888.125
$ws.Range("K122").Value = 5664.375
$ws.Range("M122").Value = -3214.375
$ws.Range("H132").Value = 7833.2144
$ws.Range("I132").Value = 7889.5835
$ws.Range("J132").Value = 7495
$ws.Range("K132").Value = 23668.7505
$ws.Range("L132").Value = 22485
$ws.Range("M132").Value = -21138.7505
$ws.Range("N132").Value = -27545

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1072.6666
$ws.Range("I22").Value = 760.1818
$ws.Range("J22").Value = 1337.0769
$ws.Range("K22").Value = 760.1818
$ws.Range("L22").Value = 1337.0769
$ws.Range("M22").Value = -465.1818
$ws.Range("N22").Value = -1927.0769
$ws.Range("H27").Value = 1072.6666
$ws.Range("I27").Value = 760.1818
$ws.Range("J27").Value = 1337.0769
$ws.Range("K27").Value = 760.1818
$ws.Range("L27").Value = 1337.0769
$ws.Range("M27").Value = -653.1818
$ws.Range("N27").Value = -1551.0769
$ws.Range("H40").Value = 5291
$ws.Range("I40").Value = 5326.1113
$ws.Range("K40").Value = 5326.1113
$ws.Range("M40").Value = -5190.1113
$ws.Range("H122").Value = 3910292.8
$ws.Range("I122").Value = 4026.5454
$ws.Range("J122").Value = 12504078
$ws.Range("K122").Value = 12079.6362
$ws.Range("L122").Value = 37512234
$ws.Range("M122").Value = -9629.636200000001
$ws.Range("N122").Value = -37517134
$ws.Range("H131").Value = 69827.86
$ws.Range("J131").Value = 76699.8
$ws.Range("L131").Value = 76699.8
$ws.Range("N131").Value = -86779.8
$ws.Range("H132").Value = 5057.625
$ws.Range("I132").Value = 2354
$ws.Range("J132").Value = 14713.429
$ws.Range("K132").Value = 7062
$ws.Range("L132").Value = 44140.287
$ws.Range("M132").Value = -4532
$ws.Range("N132").Value = -49200.287

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 2536250
$ws.Range("J124").Value = 2536250
$ws.Range("L124").Value = 2536250
$ws.Range("N124").Value = -2546070
$ws.Range("H132").Value = 6664.6665
$ws.Range("I132").Value = 6664.6665
$ws.Range("K132").Value = 19993.9995
$ws.Range("M132").Value = -17463.9995
